$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title change (appears twice: heading and bold text near the end)
Replace-Text "Play Black Bull Slot for Free - Pragmatic Play Game Review" "Play Black Bull Slot Free - Pragmatic Play Game Review"

# "What we like" bullet list
Replace-Text "Wild symbol for winning combinations" "Wild symbol and free spins increase winning opportunities"
Replace-Text "Free spins with increasing winning opportunities" "High volatility for high-value rewards"
Replace-Text "High payout potential" "Free version and bonus round options"

# "What we don't like" bullet list
Replace-Text "High volatility with infrequent rewards" "Rewards may not be frequent due to high volatility"
Replace-Text "Limited variety in symbols" "Purchasing the bonus round may require additional investment"

# Meta description (italic text)
Replace-Text "Try your luck with Black Bull online slot game from Pragmatic Play. Review includes pro and cons, volatility, free play, and purchase options." "Discover the gameplay, features, and winning potential of Black Bull slot game. Play for free!"
